$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 3.85
$ws.Range("K2").Value = 3.6
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 2.06

# Row 7
$ws.Range("F7").Value = 2.74
$ws.Range("G7").Value = 4.1
$ws.Range("I7").Value = 2.54
$ws.Range("K7").Value = 7.8
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 2.44
$ws.Range("O7").Value = 1.14
$ws.Range("R7").Value = 1.56
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000

# Row 8
$ws.Range("F8").Value = 2.1
$ws.Range("G8").Value = 2.38
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 3.55
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 4.5
$ws.Range("P8").Value = 2.4
$ws.Range("Q8").Value = 1.56

# Row 9
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 7.4
$ws.Range("H9").Value = 1.71
$ws.Range("I9").Value = 1.89
$ws.Range("J9").Value = 3.6
$ws.Range("K9").Value = 4.2
$ws.Range("P9").Value = 1.82

# Row 11
$ws.Range("P11").Value = 2.18
